# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the Notes Master  ("Office Theme" colours)
#   ppt/theme/theme2.xml -> bound to the Slide Master  ("Integral" colours)
#
# The authored edit swaps the two themes' colour schemes (and their
# display names) so that the Slide Master ends up on the stock
# "Office Theme" palette and the Notes Master ends up on the
# "Integral" palette.
#
# PowerPoint's object model only exposes theme edits through
# ThemeColorScheme(.Colors).Item(n).RGB (colour swatches) and
# ThemeFontScheme.Major/MinorFont (typeface names) - there's no
# scriptable "rename theme" property, and the Slide Master is the
# theme actually painted on every slide - so we repaint its twelve
# theme colours with the "Office Theme" palette that the swap hands
# it.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" colour scheme, in the
# standard 12-slot theme colour order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5-10 accent1..accent6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeTheme[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB long is 0x00BBGGRR
    $colors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
